$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3   = 5
    10  = 5
    15  = 2
    18  = 5
    19  = 2
    25  = 5
    30  = 5
    34  = 1
    35  = 5
    45  = 5
    51  = 5
    58  = 2
    59  = 5
    60  = 5
    64  = 2
    69  = 4
    70  = 5
    78  = 5
    80  = 5
    82  = 5
    88  = 4
    90  = 5
    94  = 5
    95  = 5
    97  = 5
    98  = 5
    101 = 5
}

foreach ($row in $changes.Keys) {
    $ws.Range("B$row").Value = $changes[$row]
}
